# Regenerate orders with updated distance/sizes.
# The distance codes D51/D64/D80 become D55/D69/D86, and the size code
# S30 becomes S31, throughout the sheet (condition labels, filenames,
# and the Distance/Size lookup columns all share these substrings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = $ws.Cells

[void]$cells.Replace("D51", "D55")
[void]$cells.Replace("D64", "D69")
[void]$cells.Replace("D80", "D86")
[void]$cells.Replace("S30", "S31")
